$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.017.17"
$ws.Range("E2").Value = "  -2.31%  "

$ws.Range("D3").Value = "2.575.52"
$ws.Range("E3").Value = "  -2.49%  "

$ws.Range("E4").Value = "  +0.07%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "517.26"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.29%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "139.07"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -4.28%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -1.70%  "

$ws.Range("D9").Value = "2.592.28"
$ws.Range("E9").Value = "  -2.34%  "

$ws.Range("E10").Value = "  -3.46%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0995"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -4.82%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.326"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -3.53%  "

$ws.Range("E13").Value = "  +0.46%  "

$ws.Range("D14").Value = "3.029.44"
$ws.Range("E14").Value = "  -2.53%  "

$ws.Range("D15").Value = "58.016.97"
$ws.Range("E15").Value = "  -2.23%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "20.08"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -4.14%  "

$ws.Range("D17").Value = "2.584.90"
$ws.Range("E17").Value = "  -2.25%  "

$ws.Range("E18").Value = "  -4.03%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "332.92"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.78%  "

$ws.Range("E20").Value = "  -3.75%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.09"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -5.07%  "

$ws.Range("E22").Value = "  -0.33%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "65.98"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  -1.57%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.996"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("E27").Value = "  -4.76%  "

$ws.Range("D28").Value = "2.688.80"
$ws.Range("E28").Value = "  -2.50%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "6.93"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -4.46%  "

$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("D31").Value = "0.0₃0716"
$ws.Range("E31").Value = "  -10.71%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.93"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -7.86%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "18.67"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -2.08%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.56"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.82%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "149.24"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("E36").Value = "  -6.71%  "

$ws.Range("E37").Value = "  -7.76%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "36.26"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.838"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -3.42%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.824"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -7.32%  "

$ws.Range("E41").Value = "  -4.73%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.51"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -4.14%  "

$ws.Range("E43").Value = "  +0.12%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "274.06"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.01%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "10.68"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.20%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.590"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.99%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.0941"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.45%  "

$ws.Range("E48").Value = "  -4.47%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "18.44"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.76%  "

$ws.Range("D50").Value = "1.974.52"
$ws.Range("E50").Value = "  -3.16%  "

$ws.Range("E51").Value = "  -5.95%  "
